$d = $word.ActiveDocument

function Insert-ItalicParagraphAfterText($anchorText, $newText) {
    # Locate the paragraph containing $anchorText, then insert a brand-new
    # paragraph right after it carrying $newText in italics.
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $para = $rng.Paragraphs(1)
    $para.Range.InsertParagraphAfter()

    $newPara = $para.Next()
    $newRng = $newPara.Range
    $newRng.Text = $newText
    $italicRng = $d.Range($newRng.Start, $newRng.Start + $newText.Length)
    $italicRng.Font.Italic = $true
}

# 1. Title change (Heading1)
$d.Content.Find.Execute("LOT2058 -  Engenharia Econômica", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LOT2058 -  Análise Técnico-Econômica de Bioprocessos", 2)

# 2. Subtitle (Heading3) change
$d.Content.Find.Execute("Fundamentals of Economic Engineering", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Engineering Economics", 2)

# 3. Ativação date change
$d.Content.Find.Execute("Ativação: 01/01/2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2025", 2)

# 6. Add a space before "NF" in the Critério paragraph
$d.Content.Find.Execute("dos exercícios individuais;NF = (0,8NP + 0,2NE)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "dos exercícios individuais; NF = (0,8NP + 0,2NE)", 2)

# 4/5/6. Insert three new italic English-translation paragraphs right after
# their Portuguese counterparts. Processed bottom-up so that earlier Find
# anchors remain valid (unaffected by later insertions further down the doc).

Insert-ItalicParagraphAfterText `
    "1.INTRODUÇÃO À ENGENHARIA ECONÔMICA: necessidades de uma análise econômica de projeto" `
    "1.INTRODUCTION TO ECONOMIC ENGINEERING: needs for an economic project analysis (economic engineering as a decision-making analysis tool);2.COST ESTIMATION: capital investment estimation (types of capital cost estimates; most common cost indices; methods for estimating capital investment); equipment cost estimation; production cost estimation;3.INTEREST: time variable (simple interest; compound interest); effective, nominal, and equivalent rates; equivalence relationships.4.CASH FLOW: cash flow diagram; cash flow equivalence; uniform and gradient series; cash flow preparation.5.DEPRECIATION: depreciation methods;6.COMPARISON OF INVESTMENT ALTERNATIVES: profitability criteria – Equivalent Uniform Annual Value (EUAV) method; Present Value (PV) method; Internal Rate of Return (IRR) method; Modified Internal Rate of Return (MIRR) method; equipment renewal and replacement; payback period;7.Project break-even point;8.Spreadsheet applications in Economic Engineering analysis."

Insert-ItalicParagraphAfterText `
    "Introdução à Engenharia Econômica; estimativa de custos; juros; fluxo de caixa" `
    "Introduction to Economic Engineering; cost estimation; interest; cash flow; depreciation; comparison of investment alternatives; break-even point; spreadsheet applications in Economic Engineering."

Insert-ItalicParagraphAfterText `
    "Capacitar os alunos a utilizarem conceitos e ferramentas de matemática financeira para avaliações econômicas de projetos de engenharia" `
    "Empower students to use concepts and tools of financial mathematics for economic evaluations of engineering projects."
